# "Add controller for demo"
# Feuil1 is a small table of avatar behaviors/text/slide-image triplets.
# Row 5 column A should drive the "swap_behave" controller instead of
# "rest_open_behave", and the infra/testing slide text has a typo fix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 5 (A5): switch the behaviour controller used for this row.
$ws.Range("A5").Value = "swap_behave"

# Row 6 (B6): fix the "Secutrité" -> "Securité" typo in the slide text.
$ws.Range("B6").Value = " Digital, Testing, Securité, et Infrastructures"

# Move the active selection to A18, matching the author's cursor position.
$ws.Range("A18").Select()
